$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.573.46"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "2.263.20"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "250.14"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "0.641"
$ws.Range("E6").Value = "  +2.21%  "
$ws.Range("E7").Value = "  +4.79%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.640"
$ws.Range("E9").Value = "  -4.47%  "
$ws.Range("D10").Value = "40.11"
$ws.Range("E10").Value = "  +2.36%  "
$ws.Range("D11").Value = "0.0964"
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.30"
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").Value = "2.605.79"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").Value = "15.02"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").Value = "0.863"
$ws.Range("E16").Value = "  -2.41%  "
$ws.Range("D17").Value = "2.283.90"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "42.522.02"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").Value = "0.0₃0989"
$ws.Range("E19").Value = "  -1.74%  "
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("D21").Value = "72.01"
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("D22").Value = "234.54"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("E23").Value = "  +2.17%  "
$ws.Range("D24").Value = "3.82"
$ws.Range("E24").Value = "  -3.07%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "11.29"
$ws.Range("E26").Value = "  -1.93%  "
$ws.Range("D27").Value = "2.37"
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").Value = "167.64"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "20.88"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").Value = "6.54"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("D32").Value = "0.0856"
$ws.Range("E32").Value = "  +6.18%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "0.124"
$ws.Range("E33").Value = "  -2.40%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "31.75"
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.50"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").Value = "4.74"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").Value = "0.0307"
$ws.Range("E38").Value = "  -4.11%  "
$ws.Range("D39").Value = "13.65"
$ws.Range("E39").Value = "  +7.46%  "
$ws.Range("E40").Value = "  -3.14%  "
$ws.Range("D41").Value = "5.87"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("D43").Value = "61.57"
$ws.Range("E43").Value = "  -1.63%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "106.41"
$ws.Range("E44").Value = "  +11.91%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "8.83"
$ws.Range("E45").Value = "  -4.86%  "
$ws.Range("E46").Value = "  -2.77%  "
$ws.Range("E47").Value = "  -2.15%  "
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("E49").Value = "  -2.82%  "
$ws.Range("D50").Value = "1.16"
$ws.Range("E50").Value = "  -3.02%  "
$ws.Range("E51").Value = "  -2.37%  "
